$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 7613
$ws.Range("E2").Value = 345
$ws.Range("F2").Value = 345
$ws.Range("G2").Value = 365
$ws.Range("H2").Value = 433
$ws.Range("I2").Value = 355
$ws.Range("J2").Value = 78
$ws.Range("K2").Value = 10333
$ws.Range("L2").Value = 4441
$ws.Range("M2").Value = 5892
$ws.Range("N2").Value = 5275
$ws.Range("O2").Value = 618
$ws.Range("P2").Value = 244
$ws.Range("Q2").Value = -328
$ws.Range("R2").Value = -692
$ws.Range("S2").Value = 851
$ws.Range("T2").Value = 692
$ws.Range("U2").Value = -1020
$ws.Range("V2").Value = 2722
$ws.Range("W2").Value = 4.53
$ws.Range("X2").Value = 5.69
$ws.Range("Y2").Value = 7.58
$ws.Range("Z2").Value = 4.52
$ws.Range("AA2").Value = 75.37
$ws.Range("AB2").Value = 2047.07
$ws.Range("AC2").Value = 3032
$ws.Range("AD2").Value = 29.19
$ws.Range("AE2").Value = 44601
$ws.Range("AF2").Value = 1.98
$ws.Range("AG2").Value = 0
$ws.Range("AH2").Value = 0
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = 11843029

# Row 3
$ws.Range("D3").Value = 13175
$ws.Range("E3").Value = 2118
$ws.Range("F3").Value = 2118
$ws.Range("G3").Value = 2096
$ws.Range("H3").Value = 1621
$ws.Range("I3").Value = 1544
$ws.Range("J3").Value = 76
$ws.Range("K3").Value = 17226
$ws.Range("L3").Value = 9778
$ws.Range("M3").Value = 7448
$ws.Range("N3").Value = 6749
$ws.Range("O3").Value = 699
$ws.Range("P3").Value = 256
$ws.Range("Q3").Value = 1017
$ws.Range("R3").Value = -1144
$ws.Range("S3").Value = 945
$ws.Range("T3").Value = 1211
$ws.Range("U3").Value = -194
$ws.Range("V3").Value = 3674
$ws.Range("W3").Value = 16.07
$ws.Range("X3").Value = 12.3
$ws.Range("Y3").Value = 25.69
$ws.Range("Z3").Value = 11.76
$ws.Range("AA3").Value = 131.28
$ws.Range("AB3").Value = 2525.43
$ws.Range("AC3").Value = 13041
$ws.Range("AD3").Value = 49.34
$ws.Range("AE3").Value = 57091
$ws.Range("AF3").Value = 11.27
$ws.Range("AG3").Value = 1728
$ws.Range("AH3").Value = 0.27
$ws.Range("AI3").Value = 13.22
$ws.Range("AJ3").Value = 11843029

# Row 4
$ws.Range("D4").Value = 8827
$ws.Range("E4").Value = 268
$ws.Range("F4").Value = 268
$ws.Range("G4").Value = 80
$ws.Range("H4").Value = 303
$ws.Range("I4").Value = 233
$ws.Range("J4").Value = 70
$ws.Range("K4").Value = 15970
$ws.Range("L4").Value = 8625
$ws.Range("M4").Value = 7345
$ws.Range("N4").Value = 6595
$ws.Range("O4").Value = 751
$ws.Range("P4").Value = 261
$ws.Range("Q4").Value = 4130
$ws.Range("R4").Value = -3551
$ws.Range("S4").Value = -629
$ws.Range("T4").Value = 2003
$ws.Range("U4").Value = 2127
$ws.Range("V4").Value = 3322
$ws.Range("W4").Value = 3.03
$ws.Range("X4").Value = 3.43
$ws.Range("Y4").Value = 3.5
$ws.Range("Z4").Value = 1.83
$ws.Range("AA4").Value = 117.42
$ws.Range("AB4").Value = 2476.93
$ws.Range("AC4").Value = 1970
$ws.Range("AD4").Value = 137.09
$ws.Range("AE4").Value = 55826
$ws.Range("AF4").Value = 4.84
$ws.Range("AG4").Value = 0
$ws.Range("AH4").Value = 0
$ws.Range("AI4").Value = 0
$ws.Range("AJ4").Value = 11843029

# Row 5
$ws.Range("D5").Value = 9166
$ws.Range("E5").Value = 822
$ws.Range("F5").Value = 822
$ws.Range("G5").Value = 588
$ws.Range("H5").Value = 690
$ws.Range("I5").Value = 604
$ws.Range("J5").Value = 85
$ws.Range("K5").Value = 16609
$ws.Range("L5").Value = 8685
$ws.Range("M5").Value = 7924
$ws.Range("N5").Value = 7117
$ws.Range("O5").Value = 807
$ws.Range("P5").Value = 279
$ws.Range("Q5").Value = -629
$ws.Range("R5").Value = -1529
$ws.Range("S5").Value = 1377
$ws.Range("T5").Value = 2471
$ws.Range("U5").Value = -3100
$ws.Range("V5").Value = 4762
$ws.Range("W5").Value = 8.960000000000001
$ws.Range("X5").Value = 7.52
$ws.Range("Y5").Value = 8.82
$ws.Range("Z5").Value = 4.23
$ws.Range("AA5").Value = 109.6
$ws.Range("AB5").Value = 2510.88
$ws.Range("AC5").Value = 5103
$ws.Range("AD5").Value = 110.28
$ws.Range("AE5").Value = 60358
$ws.Range("AF5").Value = 9.32
$ws.Range("AG5").Value = 471
$ws.Range("AH5").Value = 0.08
$ws.Range("AI5").Value = 9.19
$ws.Range("AJ5").Value = 11843029

# Row 6
$ws.Range("D6").Value = 10160
$ws.Range("E6").Value = 836
$ws.Range("F6").Value = 836
$ws.Range("G6").Value = 504
$ws.Range("H6").Value = 342
$ws.Range("I6").Value = 249
$ws.Range("K6").Value = 16912
$ws.Range("L6").Value = 9020
$ws.Range("M6").Value = 7893
$ws.Range("N6").Value = 7032
$ws.Range("P6").Value = 285
$ws.Range("Q6").Value = 260
$ws.Range("R6").Value = -1458
$ws.Range("S6").Value = 1404
$ws.Range("T6").Value = 1823
$ws.Range("U6").Value = -1563
$ws.Range("V6").Value = 6312
$ws.Range("W6").Value = 8.23
$ws.Range("X6").Value = 3.36
$ws.Range("Y6").Value = 3.52
$ws.Range("Z6").Value = 2.04
$ws.Range("AA6").Value = 114.28
$ws.Range("AB6").Value = 2468.74
$ws.Range("AC6").Value = 2102
$ws.Range("AD6").Value = 216.81
$ws.Range("AE6").Value = 59686
$ws.Range("AF6").Value = 7.63
$ws.Range("AG6").Value = 481
$ws.Range("AH6").Value = 0.11
$ws.Range("AI6").Value = 22.75
$ws.Range("AJ6").Value = 11843029

# Row 7
$ws.Range("D7").Value = 11083
$ws.Range("E7").Value = 937
$ws.Range("G7").Value = 752
$ws.Range("H7").Value = 566
$ws.Range("I7").Value = 442
$ws.Range("K7").Value = 19175
$ws.Range("L7").Value = 10789
$ws.Range("M7").Value = 8385
$ws.Range("N7").Value = 7400
$ws.Range("P7").Value = 286
$ws.Range("Q7").Value = 567
$ws.Range("R7").Value = -1441
$ws.Range("S7").Value = 1055
$ws.Range("T7").Value = 920
$ws.Range("U7").Value = -309
$ws.Range("W7").Value = 8.460000000000001
$ws.Range("X7").Value = 5.11
$ws.Range("Y7").Value = 6.13
$ws.Range("Z7").Value = 3.14
$ws.Range("AA7").Value = 128.67
$ws.Range("AC7").Value = 3735
$ws.Range("AD7").Value = 75.23
$ws.Range("AE7").Value = 62854
$ws.Range("AF7").Value = 4.47
$ws.Range("AG7").Value = 488
$ws.Range("AH7").Value = 0.17
$ws.Range("AI7").Value = 13.07

# Row 8
$ws.Range("D8").Value = 11989
$ws.Range("E8").Value = 1021
$ws.Range("G8").Value = 811
$ws.Range("H8").Value = 622
$ws.Range("I8").Value = 485
$ws.Range("K8").Value = 20221
$ws.Range("L8").Value = 11284
$ws.Range("M8").Value = 8937
$ws.Range("N8").Value = 7825
$ws.Range("P8").Value = 287
$ws.Range("Q8").Value = 953
$ws.Range("R8").Value = -1327
$ws.Range("S8").Value = 297
$ws.Range("T8").Value = 792
$ws.Range("U8").Value = -176
$ws.Range("W8").Value = 8.52
$ws.Range("X8").Value = 5.19
$ws.Range("Y8").Value = 6.37
$ws.Range("Z8").Value = 3.16
$ws.Range("AA8").Value = 126.27
$ws.Range("AC8").Value = 4091
$ws.Range("AD8").Value = 68.68000000000001
$ws.Range("AE8").Value = 66466
$ws.Range("AF8").Value = 4.23
$ws.Range("AG8").Value = 488
$ws.Range("AH8").Value = 0.17
$ws.Range("AI8").Value = 11.94

# Row 9
$ws.Range("D9").Value = 13122
$ws.Range("E9").Value = 1340
$ws.Range("G9").Value = 1009
$ws.Range("H9").Value = 761
$ws.Range("I9").Value = 650
$ws.Range("K9").Value = 19677
$ws.Range("L9").Value = 10049
$ws.Range("M9").Value = 9628
$ws.Range("N9").Value = 8360
$ws.Range("P9").Value = 288
$ws.Range("Q9").Value = 1180
$ws.Range("R9").Value = -896
$ws.Range("S9").Value = -45
$ws.Range("T9").Value = 770
$ws.Range("U9").Value = 1011
$ws.Range("W9").Value = 10.21
$ws.Range("X9").Value = 5.8
$ws.Range("Y9").Value = 8.039999999999999
$ws.Range("Z9").Value = 3.81
$ws.Range("AA9").Value = 104.37
$ws.Range("AC9").Value = 5493
$ws.Range("AD9").Value = 51.16
$ws.Range("AE9").Value = 71013
$ws.Range("AF9").Value = 3.96
$ws.Range("AG9").Value = 502
$ws.Range("AH9").Value = 0.18
$ws.Range("AI9").Value = 9.140000000000001
